# Add a new row of data (row 65) to the "Data" sheet, mirroring the
# most recent usage-log entry appended by the GCF_File_Usage logger.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 65

$ws.Cells.Item($row, 1).Value2  = 45786.928842592592   # TimeStamp
$ws.Cells.Item($row, 2).Value2  = 11                   # CC_Régularisations
$ws.Cells.Item($row, 3).Value2  = 6                     # DEB_Récurrent
$ws.Cells.Item($row, 4).Value2  = 366                   # DEB_Trans
$ws.Cells.Item($row, 5).Value2  = 545                   # ENC_Détails
$ws.Cells.Item($row, 6).Value2  = 518                   # ENC_Entête
$ws.Cells.Item($row, 7).Value2  = 575                   # FAC_Comptes_Clients
$ws.Cells.Item($row, 8).Value2  = 4305                  # FAC_Détails
$ws.Cells.Item($row, 9).Value2  = 575                   # FAC_Entête
$ws.Cells.Item($row, 10).Value2 = 2841                  # FAC_Projets_Détails
$ws.Cells.Item($row, 11).Value2 = 277                   # FAC_Projets_Entête
$ws.Cells.Item($row, 12).Value2 = 519                   # FAC_Sommaire_Taux
$ws.Cells.Item($row, 13).Value2 = 30                    # GL_EJ_Récurrente
$ws.Cells.Item($row, 14).Value2 = 4919                  # GL_Trans
$ws.Cells.Item($row, 15).Value2 = 6551                  # TEC_Local
